$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '27.023.68'
$ws.Range("E2").Value = '  +0.56%  '

# Row 3
$ws.Range("D3").Value = '1.678.81'
$ws.Range("E3").Value = '  +0.70%  '

# Row 4
$ws.Range("E4").Value = '  +0.07%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.90'
$ws.Range("E5").Value = '  +0.17%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.517'
$ws.Range("E6").Value = '  -2.43%  '

# Row 7
$ws.Range("E7").Value = '  +0.09%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.253'

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '21.47'
$ws.Range("E9").Value = '  +6.27%  '

# Row 10
$ws.Range("E10").Value = '  +0.85%  '

# Row 11
$ws.Range("E11").Value = '  -0.67%  '

# Row 12
$ws.Range("D12").Value = '1.919.67'
$ws.Range("E12").Value = '  +0.89%  '

# Row 13
$ws.Range("D13").Value = '1.666.52'
$ws.Range("E13").Value = '  -0.15%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.10'
$ws.Range("E14").Value = '  +0.53%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.532'
$ws.Range("E15").Value = '  +1.70%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.32'
$ws.Range("E16").Value = '  +0.47%  '

# Row 17
$ws.Range("D17").Value = '27.044.57'
$ws.Range("E17").Value = '  +0.54%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.17'
$ws.Range("E18").Value = '  +4.70%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '235.94'
$ws.Range("E19").Value = '  +1.80%  '

# Row 20
$ws.Range("D20").Value = '0.0₃0738'
$ws.Range("E20").Value = '  +0.68%  '

# Row 21
$ws.Range("E21").Value = '  +0.04%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.47'
$ws.Range("E22").Value = '  +0.29%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.26'
$ws.Range("E23").Value = '  +0.84%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.12'
$ws.Range("E24").Value = '  -4.06%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.97'
$ws.Range("E25").Value = '  +0.92%  '

# Row 26
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.24'
$ws.Range("E26").Value = '  +1.50%  '

# Row 27
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.57'
$ws.Range("E27").Value = '  +4.26%  '

# Row 28
$ws.Range("E28").Value = '  -2.38%  '

# Row 29
$ws.Range("E29").Value = '  +0.09%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0498'
$ws.Range("E30").Value = '  +0.47%  '

# Row 31
$ws.Range("E31").Value = '  -0.14%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.35'
$ws.Range("E32").Value = '  +0.44%  '

# Row 33
$ws.Range("D33").Value = '1.525.08'
$ws.Range("E33").Value = '  +3.86%  '

# Row 34
$ws.Range("E34").Value = '  +0.61%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.69'
$ws.Range("E35").Value = '  +4.13%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.41'
$ws.Range("E36").Value = '  -0.37%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.588'
$ws.Range("E37").Value = '  +2.75%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.916'
$ws.Range("E38").Value = '  +1.96%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0174'
$ws.Range("E39").Value = '  +3.24%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.03'
$ws.Range("E40").Value = '  +5.45%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.73'
$ws.Range("E41").Value = '  -1.46%  '

# Row 42
$ws.Range("E42").Value = '  +0.06%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '67.93'
$ws.Range("E43").Value = '  +3.31%  '

# Row 44
$ws.Range("E44").Value = '  -0.90%  '

# Row 45
$ws.Range("D45").Value = '1.825.67'
$ws.Range("E45").Value = '  +0.56%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.780'
$ws.Range("E46").Value = '  +0.08%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '90.22'
$ws.Range("E47").Value = '  -0.06%  '

# Row 48
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").Value = '0.0₆0105'
$ws.Range("E48").Value = '  -0.22%  '

# Row 49
$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.104'
$ws.Range("E49").Value = '  +3.40%  '

# Row 50
$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.52'
$ws.Range("E50").Value = '  -0.15%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.92'
$ws.Range("E51").Value = '  +4.48%  '
